# Apply the edit described by the diff:
#  1. Rename the worksheet (tab name) from "Gamma1F-HW20.xpc" to "Gamma1F".
#  2. Append a new data row (row 16) to the sheet, continuing the pattern
#     of rows 2-15, using the style of column A from the existing rows and
#     the existing shared string "HexGrid-60degTilt5degRes" for column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet/tab.
$ws.Name = "Gamma1F"

# 2. Add row 16 with the new averaged-intensity data, matching the style
#    used for column A (bordered/bold/centered header-like style) of the
#    preceding rows.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A16").Value = 14

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.18355117195812
$ws.Range("D16").Value = 0.5960752793557353
$ws.Range("E16").Value = 1.047915988402299
$ws.Range("F16").Value = 1.18355117195812
$ws.Range("G16").Value = 0.7943672782176859
$ws.Range("H16").Value = 1.124909215325272
$ws.Range("I16").Value = 1.092202312283082
$ws.Range("J16").Value = 0.5960752793557353
$ws.Range("K16").Value = 0.8219956338790171
$ws.Range("L16").Value = 1.002773402918569
$ws.Range("M16").Value = 0.9731702075903658
